$wb = $excel.ActiveWorkbook

# ---- Metadata sheet ----
$ws = $wb.Worksheets.Item("Metadata")

# Update URL value (row 2, col B)
$ws.Cells.Item(2,2).Value = "http://fhirfli.dev/fhir/ig/cicada/StructureDefinition/assessment-date"

# Update Date value (row 8, col B)
$ws.Cells.Item(8,2).Value = "2026-02-11T14:37:07-05:00"

# Insert a new row after "Contact" (row 10) for "Jurisdiction"
$ws.Rows.Item(11).Insert()

# Copy formatting from the row below (still a normal data row) so the new
# row matches the table's existing cell style/borders.
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

$ws.Cells.Item(11,1).Value = "Jurisdiction"
$ws.Cells.Item(11,2).Value = ""

# ---- Elements sheet (no content changes; shared-string table shift only) ----
